$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-12 (Generation 0-10) -> Fitness 7345
$ws.Range("C2:C12").Value = 7345

# Rows 13-38 (Generation 11-36) -> Fitness 7312
$ws.Range("C13:C38").Value = 7312
